$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 170.3
$ws.Range("J9").Value = 86.333336
$ws.Range("L9").Value = 86.333336
$ws.Range("N9").Value = -424.333336
$ws.Range("H33").Value = 996.36365
$ws.Range("I33").Value = 106.44444
$ws.Range("K33").Value = 106.44444
$ws.Range("M33").Value = 122.55556
$ws.Range("H40").Value = 7825.8276
$ws.Range("I40").Value = 3748.6667
$ws.Range("K40").Value = 3748.6667
$ws.Range("M40").Value = -3573.6667
$ws.Range("H69").Value = 6347.6665
$ws.Range("I69").Value = 1013
$ws.Range("J69").Value = 9015
$ws.Range("K69").Value = 3039
$ws.Range("L69").Value = 27045
$ws.Range("M69").Value = -2165
$ws.Range("N69").Value = -28793
$ws.Range("H72").Value = 6347.6665
$ws.Range("I72").Value = 1013
$ws.Range("J72").Value = 9015
$ws.Range("K72").Value = 9117
$ws.Range("L72").Value = 81135
$ws.Range("M72").Value = -4749
$ws.Range("N72").Value = -89871
$ws.Range("H88").Value = 2379.4614
$ws.Range("I88").Value = 5166.6665
$ws.Range("J88").Value = 1543.3
$ws.Range("K88").Value = 5166.6665
$ws.Range("L88").Value = 1543.3
$ws.Range("M88").Value = -4760.6665
$ws.Range("N88").Value = -2355.3
$ws.Range("H91").Value = 2379.4614
$ws.Range("I91").Value = 5166.6665
$ws.Range("J91").Value = 1543.3
$ws.Range("K91").Value = 5166.6665
$ws.Range("L91").Value = 1543.3
$ws.Range("M91").Value = -3762.6665
$ws.Range("N91").Value = -4351.3
$ws.Range("H98").Value = 2142.375
$ws.Range("I98").Value = 2272.1428
$ws.Range("K98").Value = 2272.1428
$ws.Range("M98").Value = -774.1428000000001
$ws.Range("H106").Value = 2963
$ws.Range("I106").Value = 2963
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2963
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2332
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 2142.375
$ws.Range("I122").Value = 2272.1428
$ws.Range("K122").Value = 6816.428400000001
$ws.Range("M122").Value = -4366.428400000001
$ws.Range("H132").Value = 2174.5
$ws.Range("I132").Value = 2115.65
$ws.Range("K132").Value = 6346.950000000001
$ws.Range("M132").Value = -3816.950000000001
$ws.Range("H135").Value = 2038.8
$ws.Range("I135").Value = 2364.6667
$ws.Range("K135").Value = 21282.0003
$ws.Range("M135").Value = -18747.0003
$ws.Range("H137").Value = 2680
$ws.Range("I137").Value = 2133.3333
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 6399.999899999999
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -3849.999899999999
$ws.Range("N137").Value = -15600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3383.3333
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10230
$ws.Range("H10").Value = 12000
$ws.Range("I10").Value = 12000
$ws.Range("K10").Value = 12000
$ws.Range("M10").Value = -11830
$ws.Range("H32").Value = 2012.3636
$ws.Range("I32").Value = 2015.2222
$ws.Range("K32").Value = 2015.2222
$ws.Range("M32").Value = -1728.2222
$ws.Range("H39").Value = 1250
$ws.Range("I39").Value = 1250
$ws.Range("K39").Value = 1250
$ws.Range("M39").Value = -730
$ws.Range("H61").Value = 1642.2858
$ws.Range("I61").Value = 1499.6666
$ws.Range("J61").Value = 2498
$ws.Range("K61").Value = 1499.6666
$ws.Range("L61").Value = 2498
$ws.Range("M61").Value = -1287.6666
$ws.Range("N61").Value = -2922
$ws.Range("H132").Value = 1356.4
$ws.Range("I132").Value = 1373.8889
$ws.Range("J132").Value = 1199
$ws.Range("K132").Value = 4121.6667
$ws.Range("L132").Value = 3597
$ws.Range("M132").Value = -1591.6667
$ws.Range("N132").Value = -8657
$ws.Range("H133").Value = 99916.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99916.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99916.8
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -104976.8
$ws.Range("H136").Value = 1642.2858
$ws.Range("I136").Value = 1499.6666
$ws.Range("J136").Value = 2498
$ws.Range("K136").Value = 4498.9998
$ws.Range("L136").Value = 7494
$ws.Range("M136").Value = -1948.9998
$ws.Range("N136").Value = -12594

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4125
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253
$ws.Range("H107").Value = 3350.2
$ws.Range("I107").Value = 1687.75
$ws.Range("K107").Value = 1687.75
$ws.Range("M107").Value = 232.25
$ws.Range("H126").Value = 89994
$ws.Range("J126").Value = 89994
$ws.Range("L126").Value = 89994
$ws.Range("N126").Value = -99874

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8433000
$ws.Range("I6").Value = 5006200
$ws.Range("J6").Value = 17000000
$ws.Range("K6").Value = 5006200
$ws.Range("L6").Value = 17000000
$ws.Range("M6").Value = -5006087
$ws.Range("N6").Value = -17000226
$ws.Range("H13").Value = 2250
$ws.Range("I13").Value = 2250
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2250
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -2111
$ws.Range("N13").ClearContents()
$ws.Range("H22").Value = 250
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 250
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -950
$ws.Range("H35").Value = 1137
$ws.Range("I35").Value = 1137
$ws.Range("K35").Value = 1137
$ws.Range("M35").Value = -843
$ws.Range("H105").Value = 6499.7144
$ws.Range("I105").Value = 1374.5
$ws.Range("K105").Value = 1374.5
$ws.Range("M105").Value = 372.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1278.3846
$ws.Range("I134").Value = 1057.3636
$ws.Range("K134").Value = 3172.0908
$ws.Range("M134").Value = -637.0907999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 190
$ws.Range("I40").Value = 75.14286
$ws.Range("K40").Value = 300.57144
$ws.Range("M40").Value = -231.57144
$ws.Range("H46").Value = 1025
$ws.Range("I46").Value = 1025
$ws.Range("K46").Value = 3075
$ws.Range("M46").Value = -2984
$ws.Range("H60").Value = 914.55554
$ws.Range("I60").Value = 531.6667
$ws.Range("J60").Value = 1680.3334
$ws.Range("K60").Value = 1595.0001
$ws.Range("L60").Value = 5041.0002
$ws.Range("M60").Value = -1344.0001
$ws.Range("N60").Value = -5543.0002
$ws.Range("H129").Value = 2737.5
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000
$ws.Range("H131").Value = 2287.3333
$ws.Range("I131").Value = 999
$ws.Range("J131").Value = 2448.375
$ws.Range("K131").Value = 2997
$ws.Range("L131").Value = 7345.125
$ws.Range("M131").Value = 2043
$ws.Range("N131").Value = -17425.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.583336
$ws.Range("J2").Value = 41.142857
$ws.Range("L2").Value = 41.142857
$ws.Range("N2").Value = -267.142857
$ws.Range("H3").Value = 278433.56
$ws.Range("I3").Value = 417267.16
$ws.Range("J3").Value = 766.3333
$ws.Range("K3").Value = 417267.16
$ws.Range("L3").Value = 766.3333
$ws.Range("M3").Value = -417151.16
$ws.Range("N3").Value = -998.3333
$ws.Range("H11").Value = 6287063
$ws.Range("I11").Value = 9364591
$ws.Range("J11").Value = 1450947.1
$ws.Range("K11").Value = 9364591
$ws.Range("L11").Value = 1450947.1
$ws.Range("M11").Value = -9364452
$ws.Range("N11").Value = -1451225.1
$ws.Range("H13").Value = 1583.3334
$ws.Range("J13").Value = 1700
$ws.Range("L13").Value = 1700
$ws.Range("N13").Value = -1978
$ws.Range("H23").Value = 3106.75
$ws.Range("J23").Value = 3106.75
$ws.Range("L23").Value = 3106.75
$ws.Range("N23").Value = -3552.75
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7587.4585
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4272
$ws.Range("H42").Value = 20025
$ws.Range("I42").Value = 20025
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 20025
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -19462
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 20025
$ws.Range("I49").Value = 20025
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 20025
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -19878
$ws.Range("N49").ClearContents()
$ws.Range("H61").Value = 3466.6667
$ws.Range("I61").Value = 3466.6667
$ws.Range("K61").Value = 3466.6667
$ws.Range("M61").Value = -3264.6667
$ws.Range("H113").Value = 3466.6667
$ws.Range("I113").Value = 3466.6667
$ws.Range("K113").Value = 3466.6667
$ws.Range("M113").Value = -1296.6667
$ws.Range("H136").Value = 4397.421
$ws.Range("I136").Value = 4159.4375
$ws.Range("J136").Value = 5666.6665
$ws.Range("K136").Value = 12478.3125
$ws.Range("L136").Value = 16999.9995
$ws.Range("M136").Value = -9928.3125
$ws.Range("N136").Value = -22099.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 17998
$ws.Range("I54").Value = 17998
$ws.Range("K54").Value = 17998
$ws.Range("M54").Value = -17478

Write-Output "Applied Kraken_Profits updates across all sheets."
